$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the key used for encrypt/decrypt test cases: aC2X0 -> aC2@X0
$ws.Range("B2").Value = "e -f sample.txt -k aC2@X0"
$ws.Range("B4").Value = "d -k aC2@X0 -f encrypted.txt"

# Update the active selection to reflect the last edited cell
$ws.Range("B4").Select()
